$d = $word.ActiveDocument

# 1) Body text: the TFS build/changeset number "C36654" -> "C36691"
#    (the "654" lives in its own bold run right after a "C36" run).
$bodyRange = $d.Content
$bodyRange.Find.ClearFormatting()
$bodyRange.Find.Execute("654", $true, $false, $false, $false, $false, $true, 1, $false, "691", 2)

# 2) Footer of the third section contains a cached PAGE field result "6"
#    that needs to read "1" after the edit.
$footerRange = $d.Sections(3).Footers(1).Range
$footerRange.Find.ClearFormatting()
$footerRange.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)
